$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2-dose series")

# Dose 1: remove latestRecAge cap (was "15 months") -> "n/a"
$ws.Range("E8").Value = "n/a"

# Dose 2: absMinAge 15 months -> 13 months (4 weeks after 12-month dose 1)
$ws.Range("B15").Value = "13 months"

# Dose 2: remove latestRecAge cap (was "6 years") -> "n/a"
$ws.Range("E15").Value = "n/a"

# Dose 2: earliestRecInt 3 months -> 4 weeks (3 months was MMRV-specific, not MMR)
$ws.Range("H16").Value = "4 weeks"
